# Adds rows 77-90 to sheet "3 nodos" and converts A76/B76 from text to numeric,
# matching "tabla de costos lista para sacar particiones".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3 nodos")

# Row 76
$ws.Range("A76").Value = 11111
$ws.Range("B76").Value = 111111
$ws.Range("C76").Value = 'Geometric'
$ws.Range("D76").Value = 2.3125
$ws.Range("E76").Value = 0.0504150390625
$ws.Range("F76").Value = '⎛ C ⎞⎛  ∅  ⎞
⎝ ∅ ⎠⎝ a,b ⎠
'

# Row 77
$ws.Range("A77").Value = 11111
$ws.Range("B77").Value = 111111
$ws.Range("C77").Value = 'Geometric'
$ws.Range("D77").Value = 2.3125
$ws.Range("E77").Value = 0.0827491283416748
$ws.Range("F77").Value = '⎛ C ⎞⎛  ∅  ⎞
⎝ ∅ ⎠⎝ a,b ⎠
'

# Row 78
$ws.Range("A78").Value = 11111
$ws.Range("B78").Value = 111111
$ws.Range("C78").Value = 'Geometric'
$ws.Range("D78").Value = 2.3125
$ws.Range("E78").Value = 0.06886148452758789
$ws.Range("F78").Value = '⎛ C ⎞⎛  ∅  ⎞
⎝ ∅ ⎠⎝ a,b ⎠
'

# Row 79
$ws.Range("A79").Value = 111111111111111
$ws.Range("B79").Value = 111111111111111
$ws.Range("C79").Value = 'Geometric'
$ws.Range("D79").Value = 14.74831676483154
$ws.Range("E79").Value = 0.5017726421356201
$ws.Range("F79").Value = '⎛ B ⎞⎛  ∅  ⎞
⎝ ∅ ⎠⎝ a,c ⎠
'

# Row 80
$ws.Range("A80").Value = 111111111111111
$ws.Range("B80").Value = 111111111111110
$ws.Range("C80").Value = 'Geometric'
$ws.Range("D80").Value = 14.73581600189209
$ws.Range("E80").Value = 0.5791902542114258
$ws.Range("F80").Value = '⎛ B ⎞⎛  ∅  ⎞
⎝ ∅ ⎠⎝ a,c ⎠
'

# Row 81
$ws.Range("A81").Value = 111111111111111
$ws.Range("B81").Value = 111111111111110
$ws.Range("C81").Value = 'Geometric'
$ws.Range("D81").Value = 14.73581600189209
$ws.Range("E81").Value = 0.5802597999572754
$ws.Range("F81").Value = '⎛ B ⎞⎛  ∅  ⎞
⎝ ∅ ⎠⎝ a,c ⎠
'

# Row 82
$ws.Range("A82").Value = 111111111111111
$ws.Range("B82").Value = 111111111111110
$ws.Range("C82").Value = 'Geometric'
$ws.Range("D82").Value = 14.73581600189209
$ws.Range("E82").Value = 0.5879178047180176
$ws.Range("F82").Value = '⎛ B ⎞⎛  ∅  ⎞
⎝ ∅ ⎠⎝ a,c ⎠
'

# Row 83
$ws.Range("A83").Value = 111111111111111
$ws.Range("B83").Value = 111111111111110
$ws.Range("C83").Value = 'Geometric'
$ws.Range("D83").Value = 14.73581600189209
$ws.Range("E83").Value = 0.470977783203125
$ws.Range("F83").Value = '⎛ B ⎞⎛  ∅  ⎞
⎝ ∅ ⎠⎝ a,c ⎠
'

# Row 84
$ws.Range("A84").Value = 111111111111111
$ws.Range("B84").Value = 111111111111111
$ws.Range("C84").Value = 'Geometric'
$ws.Range("D84").Value = 14.74831676483154
$ws.Range("E84").Value = 0.4847466945648193
$ws.Range("F84").Value = '⎛ B ⎞⎛  ∅  ⎞
⎝ ∅ ⎠⎝ a,c ⎠
'

# Row 85
$ws.Range("A85").Value = 111111111111100
$ws.Range("B85").Value = 111111111111111
$ws.Range("C85").Value = 'Geometric'
$ws.Range("D85").Value = 12.77331638336182
$ws.Range("E85").Value = 0.5350770950317383
$ws.Range("F85").Value = '⎛ B ⎞⎛  ∅  ⎞
⎝ ∅ ⎠⎝ a,c ⎠
'

# Row 86
$ws.Range("A86").Value = 111111111000000
$ws.Range("B86").Value = 1111111111
$ws.Range("C86").Value = 'Geometric'
$ws.Range("D86").Value = 8.783961296081543
$ws.Range("E86").Value = 89.12561917304993
$ws.Range("F86").Value = '⎛ B ⎞⎛  ∅  ⎞
⎝ ∅ ⎠⎝ a,c ⎠
'

# Row 87
$ws.Range("A87").Value = 111111111111111
$ws.Range("B87").Value = 111111111111111
$ws.Range("C87").Value = 'Geometric'
$ws.Range("D87").Value = 14.74831676483154
$ws.Range("E87").Value = 0.4873776435852051
$ws.Range("F87").Value = '⎛ B ⎞⎛  ∅  ⎞
⎝ ∅ ⎠⎝ a,c ⎠
'

# Row 88
$ws.Range("A88").Value = 111111111111111
$ws.Range("B88").Value = 110110110110110
$ws.Range("C88").Value = 'Geometric'
$ws.Range("D88").Value = 14.61851978302002
$ws.Range("E88").Value = 1.16145920753479
$ws.Range("F88").Value = '⎛ B ⎞⎛  ∅  ⎞
⎝ ∅ ⎠⎝ a,c ⎠
'

# Row 89
$ws.Range("A89").Value = 111111111111111
$ws.Range("B89").Value = 111111111111111
$ws.Range("C89").Value = 'Geometric'
$ws.Range("D89").Value = 14.74831676483154
$ws.Range("E89").Value = 0.4490947723388672
$ws.Range("F89").Value = '⎛ B ⎞⎛  ∅  ⎞
⎝ ∅ ⎠⎝ a,c ⎠
'

# Row 90
$ws.Range("A90:B90").NumberFormat = "@"
$ws.Range("A90").Value = '111111111111111'
$ws.Range("B90").Value = '110110110110110'
$ws.Range("A90:B90").ClearFormats()
$ws.Range("C90").Value = 'Geometric'
$ws.Range("D90").Value = 14.61851978302002
$ws.Range("E90").Value = 1.101146697998047
$ws.Range("F90").Value = '⎛ B ⎞⎛  ∅  ⎞
⎝ ∅ ⎠⎝ a,c ⎠
'
